# Refresh the crypto Price (D) / Volume(1h) (E) columns, and the one
# row-49/51 coin swap (RocketPoolETH -> BEAM -> SEI reshuffle), to match
# the scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (e.g. "69.03", "0.0847"), not real
# numbers. Force the whole column to keep a Text format before writing the
# refreshed values, otherwise Excel would auto-convert look-alike numerics
# and silently drop significant trailing zeros (e.g. "69.00" -> 69,
# "3.50" -> 3.5) or switch the cell's stored type from Text to Number.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.745.13"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "2.929.86"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "351.83"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("D6").Value = "107.03"
$ws.Range("E6").Value = "  -5.58%  "

$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  -1.87%  "

$ws.Range("D10").Value = "37.78"
$ws.Range("E10").Value = "  -4.50%  "

$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").Value = "0.0847"
$ws.Range("E12").Value = "  -3.65%  "

$ws.Range("D13").Value = "18.93"
$ws.Range("E13").Value = "  -5.47%  "

$ws.Range("D14").Value = "3.393.73"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "7.56"

$ws.Range("D16").Value = "2.933.41"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").Value = "0.971"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").Value = "51.726.78"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("D19").Value = "3.50"
$ws.Range("E19").Value = "  +5.90%  "

$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -3.14%  "

$ws.Range("D21").Value = "13.42"
$ws.Range("E21").Value = "  -4.89%  "

$ws.Range("E22").Value = "  -1.61%  "

$ws.Range("D23").Value = "69.00"
$ws.Range("E23").Value = "  -3.04%  "

$ws.Range("D24").Value = "262.27"
$ws.Range("E24").Value = "  -2.79%  "

$ws.Range("E25").Value = "  -3.41%  "

$ws.Range("D26").Value = "0.174"
$ws.Range("E26").Value = "  -4.24%  "

$ws.Range("D27").Value = "26.54"
$ws.Range("E27").Value = "  -1.18%  "

$ws.Range("D28").Value = "7.55"
$ws.Range("E28").Value = "  +7.60%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").Value = "10.23"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("D32").Value = "6.15"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("D33").Value = "35.62"
$ws.Range("E33").Value = "  -4.29%  "

$ws.Range("D34").Value = "2.16"
$ws.Range("E34").Value = "  -5.02%  "

$ws.Range("E35").Value = "  -4.25%  "

$ws.Range("D36").Value = "0.0425"
$ws.Range("E36").Value = "  -6.17%  "

$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("E38").Value = "  -6.76%  "

$ws.Range("E39").Value = "  -4.13%  "

$ws.Range("D40").Value = "17.68"
$ws.Range("E40").Value = "  -5.68%  "

$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").Value = "22.66"
$ws.Range("E43").Value = "  -2.30%  "

$ws.Range("D44").Value = "119.70"
$ws.Range("E44").Value = "  +1.88%  "

$ws.Range("E45").Value = "  -1.48%  "

$ws.Range("D46").Value = "2.098.98"
$ws.Range("E46").Value = "  -3.96%  "

$ws.Range("E47").Value = "  -5.85%  "

$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -7.67%  "

$ws.Range("B49").Value = "BEAM"
$ws.Range("C49").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D49").Value = "0.0345"
$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("D50").Value = "0.239"
$ws.Range("E50").Value = "  -4.79%  "

$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").Value = "0.893"
$ws.Range("E51").Value = "  -6.25%  "
